$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# New comment for the "2^n recursive subsequences" row.
$ws.Range("C20").Value = "2^n recursive solution working"

# D20 stays empty but picks up the same "Status" highlight formatting
# used by the other data rows (D3, D4, D5, D7, D16, D17). Copy/paste
# the format from an existing cell so the existing style is reused
# instead of creating a new one.
$ws.Range("D17").Copy()
$ws.Range("D20").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Move the active selection to D22, matching the authored workbook state.
$ws.Range("D22").Select()
